$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = 3317384.86
$ws.Range("C9").Value = 514326.99
$ws.Range("D9").Value = 3831711.85
$ws.Range("E9").Value = 13.42290365597298
$ws.Range("F9").Value = 86.57709634402703
$ws.Range("G9").Value = -50.29286096910227
$ws.Range("H9").Value = -40.09261484917015
$ws.Range("I9").Value = 33284
$ws.Range("J9").Value = 1399
$ws.Range("K9").Value = 34683
$ws.Range("L9").Value = 23936
$ws.Range("M9").Value = 160.0815445354278
$ws.Range("N9").Value = 9.290930951919151
